$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; existing columns A:W shift to B:X,
# carrying their values/formatting with them.
$ws.Columns("A").Insert()

# Populate the new column A with the "Match ID" header and the match id (14)
# for every data row, applying bold formatting to match the rest of the
# header/id styling.
$ws.Range("A1").Value = "Match ID"
$ws.Range("A1:A19").Font.Bold = $true

$ws.Range("A4:A19").Value = 14

# Row 20 is a hidden totals row; writing into a hidden row can perturb its
# row height, so temporarily unhide it while setting the value.
$ws.Rows(20).Hidden = $false
$ws.Range("A20").Value = 14
$ws.Rows(20).Hidden = $true

$ws.Range("A1:A19").Select() | Out-Null
